$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(34, 1).Value = "Gennaro Bullo"
$ws.Cells.Item(34, 2).Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Cells.Item(34, 3).Value = "ENRICO BORDIGNON | Pinguini Trentini"
$ws.Cells.Item(34, 4).Value = "Leonardo Viola | Shark Attack"
$ws.Cells.Item(34, 5).Value = "Randy Cobbinah | MAI UNA GIOIA"
$ws.Cells.Item(34, 6).Value = "Matteo Mazzola | MediaserT"
